$wb = $excel.ActiveWorkbook

# Sheet 1: "展览" - update "想去人数" (want-to-go count) figures
$ws = $wb.Worksheets.Item(1)
$ws.Range("F5").Value = 1208
$ws.Range("F7").Value = 4532
$ws.Range("F8").Value = 2679
$ws.Range("F9").Value = 64
$ws.Range("F10").Value = 2646
$ws.Range("F15").Value = 700
$ws.Range("F16").Value = 448
$ws.Range("F17").Value = 167
$ws.Range("F18").Value = 356
$ws.Range("F22").Value = 51
$ws.Range("F26").Value = 599
$ws.Range("F27").Value = 720
$ws.Range("F28").Value = 129
$ws.Range("F29").Value = 16
$ws.Range("F30").Value = 463
$ws.Range("F31").Value = 1641
$ws.Range("F32").Value = 1289
$ws.Range("F33").Value = 248
$ws.Range("F35").Value = 1322
$ws.Range("F36").Value = 2178
$ws.Range("F37").Value = 334
$ws.Range("F38").Value = 19
$ws.Range("F39").Value = 574
$ws.Range("F43").Value = 723
$ws.Range("F44").Value = 1404
$ws.Range("F45").Value = 160
$ws.Range("F47").Value = 459
$ws.Range("F48").Value = 64
$ws.Range("F49").Value = 90

# Sheet 2: "演出" - update "想去人数" figures
$ws = $wb.Worksheets.Item(2)
$ws.Range("F8").Value = 8
$ws.Range("F12").Value = 14

# Sheet 4: "全部类型" - update "想去人数" figures
$ws = $wb.Worksheets.Item(4)
$ws.Range("F5").Value = 4532
$ws.Range("F6").Value = 2679
$ws.Range("F7").Value = 2646
$ws.Range("F11").Value = 700
$ws.Range("F12").Value = 448
$ws.Range("F13").Value = 167
$ws.Range("F14").Value = 356
$ws.Range("F18").Value = 51
$ws.Range("F21").Value = 599
$ws.Range("F22").Value = 720
$ws.Range("F23").Value = 129
$ws.Range("F27").Value = 16
$ws.Range("F28").Value = 463
$ws.Range("F29").Value = 1641
$ws.Range("F30").Value = 1289
$ws.Range("F31").Value = 248
$ws.Range("F33").Value = 2178
$ws.Range("F34").Value = 334
$ws.Range("F37").Value = 19
$ws.Range("F38").Value = 14
$ws.Range("F39").Value = 574
$ws.Range("F43").Value = 723
$ws.Range("F44").Value = 1404
$ws.Range("F46").Value = 160
$ws.Range("F47").Value = 459
$ws.Range("F48").Value = 90
